$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feed items discovered by the workflow run: Thermo Fisher / Bayer lung cancer CDx story,
# picked up from both GenomeWeb and 360Dx.
$link64 = "https://www.genomeweb.com/regulatory-news-fda-approvals/thermo-fisher-receives-fda-approval-ngs-based-cdx-bayer-lung-cancer"
$link65 = "https://www.360dx.com/regulatory-news-fda-approvals/thermo-fisher-receives-fda-approval-ngs-based-cdx-bayer-lung-cancer"
$keyword = "CDx"
$title = "Thermo Fisher Receives FDA Approval for NGS-Based CDx for Bayer Lung Cancer Therapy"

# Row 64
$ws.Range("A64").Value = $link64
$ws.Hyperlinks.Add($ws.Range("A64"), $link64)
$ws.Range("A64").Style = $ws.Range("A63").Style
$ws.Range("B64").Value = $keyword
$ws.Range("C64").Value = $title

# Row 65
$ws.Range("A65").Value = $link65
$ws.Hyperlinks.Add($ws.Range("A65"), $link65)
$ws.Range("A65").Style = $ws.Range("A63").Style
$ws.Range("B65").Value = $keyword
$ws.Range("C65").Value = $title
